$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new entry row 77: date, hours, task description
$ws.Cells.Item(77, 1).Value = 45401
$ws.Cells.Item(77, 1).NumberFormat = "d\-mmm"
$ws.Cells.Item(77, 2).Value = 8
$ws.Cells.Item(77, 3).Value = "RQ4 anfangen -> QSE QB geht nicht"

# Update selection to match author's final cursor position
$null = $ws.Range("G70").Select()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1

$null = $wb.Save()
